$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so numeric-looking strings (e.g. "1.002", "22.353.68")
# are written as text, matching the source data (t="inlineStr" in the XML),
# then restore the default "Normal" style so no stray style index is left on the cells.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '22.353.68'
$ws.Range("E2").Value = '  -4.93%  '
$ws.Range("D3").Value = '1.562.38'
$ws.Range("E3").Value = '  -5.15%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '289.61'
$ws.Range("E6").Value = '  -3.56%  '
$ws.Range("D7").Value = '0.3713'
$ws.Range("E7").Value = '  -2.27%  '
$ws.Range("D8").Value = '49.09'
$ws.Range("E8").Value = '  -2.68%  '
$ws.Range("D9").Value = '0.3397'
$ws.Range("E9").Value = '  -3.29%  '
$ws.Range("D10").Value = '1.166'
$ws.Range("E10").Value = '  -4.59%  '
$ws.Range("D11").Value = '0.07633'
$ws.Range("E11").Value = '  -5.43%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '21.39'
$ws.Range("E13").Value = '  -3.22%  '
$ws.Range("D14").Value = '6.044'
$ws.Range("E14").Value = '  -4.31%  '
$ws.Range("D15").Value = '6.913'
$ws.Range("E15").Value = '  -4.77%  '
$ws.Range("D16").Value = '1.563.74'
$ws.Range("E16").Value = '  -4.62%  '
$ws.Range("D17").Value = '0.00001127'
$ws.Range("E17").Value = '  -7.20%  '
$ws.Range("D18").Value = '89.99'
$ws.Range("E18").Value = '  -5.26%  '
$ws.Range("D19").Value = '0.06726'
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").Value = '6.228'
$ws.Range("E21").Value = '  -5.98%  '
$ws.Range("D22").Value = '16.53'
$ws.Range("E22").Value = '  -5.39%  '
$ws.Range("D23").Value = '0.5306'
$ws.Range("E23").Value = '  -7.18%  '
$ws.Range("D24").Value = '12.00'
$ws.Range("E24").Value = '  -3.83%  '
$ws.Range("D25").Value = '22.353.10'
$ws.Range("E25").Value = '  -4.95%  '
$ws.Range("D26").Value = '2.398'
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("D27").Value = '2.822'
$ws.Range("E27").Value = '  -5.10%  '
$ws.Range("D28").Value = '20.17'
$ws.Range("E28").Value = '  -3.92%  '
$ws.Range("D29").Value = '145.44'
$ws.Range("E29").Value = '  -3.76%  '
$ws.Range("D30").Value = '4.988'
$ws.Range("E30").Value = '  -3.67%  '
$ws.Range("D31").Value = '125.30'
$ws.Range("E31").Value = '  -4.82%  '
$ws.Range("D32").Value = '1.739.36'
$ws.Range("E32").Value = '  -5.11%  '
$ws.Range("D33").Value = '6.192'
$ws.Range("E33").Value = '  -9.55%  '
$ws.Range("D34").Value = '2.006'
$ws.Range("E34").Value = '  -6.08%  '
$ws.Range("D35").Value = '1.002'
$ws.Range("E35").Value = '  +1.19%  '
$ws.Range("D36").Value = '10.03'
$ws.Range("E36").Value = '  -10.69%  '
$ws.Range("D37").Value = '0.08438'
$ws.Range("E37").Value = '  -3.98%  '
$ws.Range("D38").Value = '0.02531'
$ws.Range("E38").Value = '  -6.35%  '
$ws.Range("D39").Value = '0.2322'
$ws.Range("E39").Value = '  -4.21%  '
$ws.Range("D40").Value = '5.519'
$ws.Range("E40").Value = '  -6.85%  '
$ws.Range("D41").Value = '0.06384'
$ws.Range("E41").Value = '  -6.05%  '
$ws.Range("D42").Value = '1.297'
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").Value = '11.69'
$ws.Range("E43").Value = '  -8.92%  '
$ws.Range("D44").Value = '0.6343'
$ws.Range("E44").Value = '  -7.66%  '
$ws.Range("D45").Value = '14.15'
$ws.Range("E45").Value = '  -8.69%  '
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").Value = '0.5967'
$ws.Range("E47").Value = '  -6.56%  '
$ws.Range("D48").Value = '3.756'
$ws.Range("E48").Value = '  -4.24%  '
$ws.Range("D49").Value = '2.092'
$ws.Range("E49").Value = '  -7.01%  '
$ws.Range("D50").Value = '1.268'
$ws.Range("E50").Value = '  +2.89%  '
$ws.Range("D51").Value = '124.60'
$ws.Range("E51").Value = '  -2.06%  '

$rng.Style = "Normal"
